# Warrant Modification Request Mapping.xlsx
# Commit: "Added Vehicle Registartion State Code to 'issued' and 'modification' SSPs."
#
# The sheet lists data-element mappings as one table row per element.
# A new mapping row ("Vehicle Registration State" -> MI -> NCIC/LIS jurisdiction
# code xpath) is inserted right above the existing row 46
# ("Vehicle Primary Color Code Text"), shifting every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 46 (pushes old row 46.. down to 47..)
$ws.Rows("46:46").Insert()

# Populate the new row's three data cells (Class/Element column is blank,
# same as its neighbouring "Vehicle ..." rows).
$ws.Range("C46").Value = "Vehicle Registration State"
$ws.Range("E46").Value = "MI"
$ws.Range("F46").Value = "wm-req-doc:WarrantModificationRequest/j:ConveyanceRegistration[@structures:id=/wm-req-doc:WarrantModificationRequest/j:ConveyanceRegistrationAssociation/j:ItemRegistration/@structures:ref]/j:JurisdictionNCICLISCode"

# Match the row's taller wrapped-text height seen in the authored sheet.
$ws.Rows("46:46").RowHeight = 56

# Move the frozen-pane view / selection to where the author was working
# when they saved (near the newly inserted row).
$ws.Range("E46").Select() | Out-Null
